$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xgb = "XGB - Extreme Boosting Classifier"
$rf  = "Random Forest"
$svm = "SVM - Support Vector Machine"

$standardScaler   = "StandardScaler"
$robustScaler     = "RobustScaler"
$powerTransformer = "PowerTransformer"
$minMaxScaler     = "MinMaxScaler"
$normalizer       = "Normalizer"

# Full results table (rows 2-16): name / standardization / f1 score
$colA = @($xgb, $xgb, $xgb, $xgb, $rf, $rf, $rf, $rf, $svm, $svm, $svm, $svm, $rf, $xgb, $svm)
$colB = @($standardScaler, $robustScaler, $powerTransformer, $minMaxScaler, $standardScaler, $robustScaler, $powerTransformer, $minMaxScaler, $standardScaler, $robustScaler, $minMaxScaler, $powerTransformer, $normalizer, $normalizer, $normalizer)
$colC = @(83.44, 83.44, 83.44, 83.44, 82.72, 82.42, 81.25, 80.77, 76.92, 76.02, 75.29000000000001, 73.05, 72.19, 71.01000000000001, 70.53)

# Write column by column (algorithm names, then scaler names, then scores)
for ($i = 0; $i -lt $colA.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt $colB.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $colB[$i]
}
for ($i = 0; $i -lt $colC.Count; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $colC[$i]
}
